$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "escaped @" test rows currently sitting right before the
# trailing "+38(099) ..." anchor row (rows 39-41) need to move up to
# right after the other "foo.bar@gmail..." rows (becoming rows 24-26),
# and every pattern that uses the old backslash escape ("\@") must be
# rewritten to use the bracket escape ("[@]") instead.

# Switch the escape syntax on every cell that uses it, in place, before
# moving anything around.
$ws.Range("D41").Value = "@+[@]@+.@+"
$ws.Range("D22").Value = "?+[@]?+.?+"
$ws.Range("D23").Value = "?+[@]?+.?+"
$ws.Range("D40").Value = "?+[@]?+.?+"
$ws.Range("D39").Value = "*[@]*.*"

# Capture the values of the three rows that will be relocated.
$movedC1 = $ws.Range("C39").Value2
$movedD1 = $ws.Range("D39").Value2
$movedE1 = $ws.Range("E39").Value2

$movedC2 = $ws.Range("C40").Value2
$movedD2 = $ws.Range("D40").Value2
$movedE2 = $ws.Range("E40").Value2

$movedC3 = $ws.Range("C41").Value2
$movedD3 = $ws.Range("D41").Value2
$movedE3 = $ws.Range("E41").Value2

# Remove the three rows from their old location (rows below shift up).
$ws.Rows("39:41").Delete()

# Make room for the three rows at their new location (rows from here
# down shift down by three).
$ws.Rows("24:26").Insert()

# Write the relocated rows back into their new home.
$ws.Range("C24").Value = $movedC1
$ws.Range("D24").Value = $movedD1
$ws.Range("E24").Value = $movedE1

$ws.Range("C25").Value = $movedC2
$ws.Range("D25").Value = $movedD2
$ws.Range("E25").Value = $movedE2

$ws.Range("C26").Value = $movedC3
$ws.Range("D26").Value = $movedD3
$ws.Range("E26").Value = $movedE3

# Update the remembered selection.
$ws.Range("L30").Select()
